$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.920.42"
$ws.Range("E2").Value = "  +2.21%  "

$ws.Range("D3").Value = "2.059.48"
$ws.Range("E3").Value = "  +2.02%  "

$ws.Range("E4").Value = "  -0.32%  "

$ws.Range("D5").Value = "'230.75"
$ws.Range("E5").Value = "  +1.82%  "

$ws.Range("D6").Value = "'0.616"
$ws.Range("E6").Value = "  +1.41%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'57.37"
$ws.Range("E8").Value = "  +4.85%  "

$ws.Range("D9").Value = "'0.386"
$ws.Range("E9").Value = "  +1.67%  "

$ws.Range("D10").Value = "'0.0804"
$ws.Range("E10").Value = "  +2.55%  "

$ws.Range("E11").Value = "  -0.67%  "

$ws.Range("D12").Value = "2.364.78"
$ws.Range("E12").Value = "  +2.14%  "

$ws.Range("D13").Value = "'14.59"
$ws.Range("E13").Value = "  +2.13%  "

$ws.Range("D14").Value = "'20.66"
$ws.Range("E14").Value = "  +1.37%  "

$ws.Range("D15").Value = "'5.27"
$ws.Range("E15").Value = "  +2.56%  "

$ws.Range("D16").Value = "'0.749"
$ws.Range("E16").Value = "  +1.10%  "

$ws.Range("D17").Value = "2.058.22"
$ws.Range("E17").Value = "  +1.90%  "

$ws.Range("D18").Value = "37.874.75"
$ws.Range("E18").Value = "  +2.35%  "

$ws.Range("D19").Value = "'6.25"
$ws.Range("E19").Value = "  +1.85%  "

$ws.Range("D20").Value = "'69.69"
$ws.Range("E20").Value = "  +1.20%  "

$ws.Range("D21").Value = "0.0₃0830"
$ws.Range("E21").Value = "  +1.17%  "

$ws.Range("D22").Value = "'224.65"
$ws.Range("E22").Value = "  -0.68%  "

$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("E24").Value = "  +1.52%  "

$ws.Range("E25").Value = "  +4.15%  "

$ws.Range("D26").Value = "'9.28"
$ws.Range("E26").Value = "  +0.71%  "

$ws.Range("D27").Value = "'165.81"
$ws.Range("E27").Value = "  +0.41%  "

$ws.Range("E28").Value = "  +7.02%  "

$ws.Range("D29").Value = "'19.12"
$ws.Range("E29").Value = "  +2.10%  "

$ws.Range("D30").Value = "'1.36"
$ws.Range("E30").Value = "  +1.42%  "

$ws.Range("E31").Value = "  +1.28%  "

$ws.Range("D32").Value = "'4.54"
$ws.Range("E32").Value = "  +0.64%  "

$ws.Range("D33").Value = "'0.0614"
$ws.Range("E33").Value = "  -0.16%  "

$ws.Range("D34").Value = "'4.56"
$ws.Range("E34").Value = "  +2.86%  "

$ws.Range("E35").Value = "  +9.67%  "

$ws.Range("D36").Value = "'2.37"
$ws.Range("E36").Value = "  +0.93%  "

$ws.Range("E37").Value = "  +11.10%  "

$ws.Range("D38").Value = "'3.30"
$ws.Range("E38").Value = "  +4.90%  "

$ws.Range("E39").Value = "  -0.11%  "

$ws.Range("D40").Value = "'99.32"
$ws.Range("E40").Value = "  +4.24%  "

$ws.Range("B41").Value = "Cronos"
$ws.Range("C41").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D41").Value = "'0.0956"
$ws.Range("E41").Value = "  +3.77%  "

$ws.Range("D42").Value = "'0.0218"
$ws.Range("E42").Value = "  +0.09%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.482.85"
$ws.Range("E43").Value = "  +0.37%  "

$ws.Range("D44").Value = "'16.76"
$ws.Range("E44").Value = "  +0.57%  "

$ws.Range("E45").Value = "  +2.19%  "

$ws.Range("E46").Value = "  -0.59%  "

$ws.Range("E47").Value = "  +13.64%  "

$ws.Range("E48").Value = "  +1.55%  "

$ws.Range("D49").Value = "'7.14"
$ws.Range("E49").Value = "  -1.64%  "

$ws.Range("D50").Value = "'2.95"
$ws.Range("E50").Value = "  +0.69%  "

$ws.Range("D51").Value = "2.246.07"
$ws.Range("E51").Value = "  +1.94%  "
